$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.798.27'
$ws.Range("E2").Value = '  +2.27%  '
$ws.Range("D3").Value = '2.083.69'
$ws.Range("E3").Value = '  +4.38%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.619'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.50'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.54%  '
$ws.Range("E9").Value = '  +3.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.28'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.01%  '
$ws.Range("E11").Value = '  +1.67%  '
$ws.Range("E12").Value = '  +3.83%  '
$ws.Range("D13").Value = '2.392.85'
$ws.Range("E13").Value = '  +4.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.49'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.11'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.784'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.33%  '
$ws.Range("E17").Value = '  +4.34%  '
$ws.Range("D18").Value = '2.113.61'
$ws.Range("E18").Value = '  +4.69%  '
$ws.Range("D19").Value = '37.903.86'
$ws.Range("E19").Value = '  +2.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +20.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '68.97'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.53%  '
$ws.Range("D22").Value = '0.0₃0818'
$ws.Range("E22").Value = '  +0.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '225.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.31%  '
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.74%  '
$ws.Range("E28").Value = '  +2.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.132'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.42'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.39'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.32%  '
$ws.Range("E32").Value = '  +1.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.50'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.42%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0632'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.16%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.61'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +12.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.49'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.72%  '
$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.97'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.95%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.36'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.86%  '
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("E41").Value = '  -3.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0968'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.96%  '
$ws.Range("D43").Value = '1.486.44'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '95.85'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.52%  '
$ws.Range("E45").Value = '  +4.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +26.56%  '
$ws.Range("E48").Value = '  +0.91%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.38'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.74%  '
$ws.Range("E50").Value = '  +2.31%  '
$ws.Range("E51").Value = '  +1.78%  '
